# Auto-generated Excel COM-interop script to apply numeric updates
# to the Tonberry_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1663.3334
$ws.Range("J39").Value = 1672.8572
$ws.Range("L39").Value = 5018.571599999999
$ws.Range("N39").Value = -5610.571599999999

$ws.Range("H80").Value = 2310.875
$ws.Range("I80").Value = 3567
$ws.Range("K80").Value = 10701
$ws.Range("M80").Value = -9703

$ws.Range("H83").Value = 2310.875
$ws.Range("I83").Value = 3567
$ws.Range("K83").Value = 32103
$ws.Range("M83").Value = -27111

$ws.Range("H98").Value = 1458.6296
$ws.Range("I98").Value = 1235.32
$ws.Range("K98").Value = 1235.32
$ws.Range("M98").Value = 262.6800000000001

$ws.Range("H116").Value = 11094.091
$ws.Range("I116").Value = 21682
$ws.Range("K116").Value = 21682
$ws.Range("M116").Value = -18240

$ws.Range("H122").Value = 1458.6296
$ws.Range("I122").Value = 1235.32
$ws.Range("K122").Value = 3705.96
$ws.Range("M122").Value = -1255.96

$ws.Range("H132").Value = 891.4902
$ws.Range("I132").Value = 828.3333
$ws.Range("J132").Value = 1902
$ws.Range("K132").Value = 2484.9999
$ws.Range("L132").Value = 5706
$ws.Range("M132").Value = 45.0001000000002
$ws.Range("N132").Value = -10766

$ws.Range("H137").Value = 1797.0526
$ws.Range("I137").Value = 1737.375
$ws.Range("J137").Value = 1840.4546
$ws.Range("K137").Value = 5212.125
$ws.Range("L137").Value = 5521.3638
$ws.Range("M137").Value = -2662.125
$ws.Range("N137").Value = -10621.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 223034.88
$ws.Range("I2").Value = 293024.9
$ws.Range("K2").Value = 293024.9
$ws.Range("M2").Value = -292911.9

$ws.Range("H4").Value = 500.4
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 502
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 502
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = -734

$ws.Range("H32").Value = 17898.771
$ws.Range("I32").Value = 14166.5
$ws.Range("J32").Value = 22309.637
$ws.Range("K32").Value = 14166.5
$ws.Range("L32").Value = 22309.637
$ws.Range("M32").Value = -13879.5
$ws.Range("N32").Value = -22883.637

$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976

$ws.Range("H55").Value = 19000
$ws.Range("J55").Value = 19000
$ws.Range("L55").Value = 19000
$ws.Range("N55").Value = -19630

$ws.Range("H74").Value = 1003.2766
$ws.Range("I74").Value = 765.9286
$ws.Range("K74").Value = 765.9286
$ws.Range("M74").Value = 108.0714

$ws.Range("H77").Value = 1003.2766
$ws.Range("I77").Value = 765.9286
$ws.Range("K77").Value = 3829.643
$ws.Range("M77").Value = 538.357

$ws.Range("H80").Value = 49333.332
$ws.Range("J80").Value = 49333.332
$ws.Range("L80").Value = 49333.332
$ws.Range("N80").Value = -51329.332

$ws.Range("H83").Value = 49333.332
$ws.Range("J83").Value = 49333.332
$ws.Range("L83").Value = 147999.996
$ws.Range("N83").Value = -157983.996

$ws.Range("H116").Value = 223034.88
$ws.Range("I116").Value = 293024.9
$ws.Range("K116").Value = 293024.9
$ws.Range("M116").Value = -290730.9

$ws.Range("H122").Value = 1902.4
$ws.Range("I122").Value = 1902.4
$ws.Range("K122").Value = 5707.200000000001
$ws.Range("M122").Value = -3257.200000000001

$ws.Range("H132").Value = 3093.9048
$ws.Range("I132").Value = 4418.75
$ws.Range("K132").Value = 13256.25
$ws.Range("M132").Value = -10726.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 223034.88
$ws.Range("I3").Value = 293024.9
$ws.Range("K3").Value = 293024.9
$ws.Range("M3").Value = -292910.9

$ws.Range("H22").Value = 787
$ws.Range("J22").Value = 846.2
$ws.Range("L22").Value = 846.2
$ws.Range("N22").Value = -1192.2

$ws.Range("H64").Value = 232.5
$ws.Range("I64").Value = 301.75
$ws.Range("J64").Value = 94
$ws.Range("K64").Value = 301.75
$ws.Range("L64").Value = 94
$ws.Range("M64").Value = -76.75
$ws.Range("N64").Value = -544

$ws.Range("H67").Value = 232.5
$ws.Range("I67").Value = 301.75
$ws.Range("J67").Value = 94
$ws.Range("K67").Value = 301.75
$ws.Range("L67").Value = 94
$ws.Range("M67").Value = 478.25
$ws.Range("N67").Value = -1654

$ws.Range("H94").Value = 1674.6666
$ws.Range("I94").Value = 752.5
$ws.Range("J94").Value = 2135.75
$ws.Range("K94").Value = 752.5
$ws.Range("L94").Value = 2135.75
$ws.Range("M94").Value = -301.5
$ws.Range("N94").Value = -3037.75

$ws.Range("H137").Value = 45793.332
$ws.Range("J137").Value = 45793.332
$ws.Range("L137").Value = 45793.332
$ws.Range("N137").Value = -55993.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3551.1304
$ws.Range("I31").Value = 1573.8182
$ws.Range("K31").Value = 1573.8182
$ws.Range("M31").Value = -1278.8182

$ws.Range("H34").Value = 3551.1304
$ws.Range("I34").Value = 1573.8182
$ws.Range("K34").Value = 1573.8182
$ws.Range("M34").Value = -1371.8182

$ws.Range("H62").Value = 3517
$ws.Range("I62").Value = 2995.9092
$ws.Range("J62").Value = 4950
$ws.Range("K62").Value = 2995.9092
$ws.Range("L62").Value = 4950
$ws.Range("M62").Value = -2371.9092
$ws.Range("N62").Value = -6198

$ws.Range("H65").Value = 3517
$ws.Range("I65").Value = 2995.9092
$ws.Range("J65").Value = 4950
$ws.Range("K65").Value = 14979.546
$ws.Range("L65").Value = 24750
$ws.Range("M65").Value = -11859.546
$ws.Range("N65").Value = -30990

$ws.Range("H99").Value = 2238.3333
$ws.Range("I99").Value = 2206.7778
$ws.Range("K99").Value = 2206.7778
$ws.Range("M99").Value = -708.7777999999998

$ws.Range("H126").Value = 2238.3333
$ws.Range("I126").Value = 2206.7778
$ws.Range("K126").Value = 6620.3334
$ws.Range("M126").Value = -4150.3334

$ws.Range("H141").Value = 66684.75
$ws.Range("J141").Value = 66201.55
$ws.Range("L141").Value = 66201.55
$ws.Range("N141").Value = -76561.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1300
$ws.Range("I68").Value = 900
$ws.Range("K68").Value = 2700
$ws.Range("M68").Value = -1889

$ws.Range("H71").Value = 1300
$ws.Range("I71").Value = 900
$ws.Range("K71").Value = 8100
$ws.Range("M71").Value = -4044

$ws.Range("H131").Value = 21552.045
$ws.Range("I131").Value = 493.75
$ws.Range("J131").Value = 23606.512
$ws.Range("K131").Value = 1481.25
$ws.Range("L131").Value = 70819.53599999999
$ws.Range("M131").Value = 3558.75
$ws.Range("N131").Value = -80899.53599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 15837.375
$ws.Range("J136").Value = 15837.375
$ws.Range("L136").Value = 47512.125
$ws.Range("N136").Value = -52612.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1524.7179
$ws.Range("I132").Value = 1278.6875
$ws.Range("K132").Value = 3836.0625
$ws.Range("M132").Value = -1306.0625
